$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two blank helper columns (K and L) sitting before the
# "hg18 build coordinates" column (originally column M). Remove those two
# empty columns so the coordinate column shifts left into column K and
# the sheet's used range shrinks from A1:P50 to A1:N50.
$ws.Range("K1:L1").EntireColumn.Delete()

# Leave the selection where the author left it after making the edit.
$null = $ws.Range("L12").Select()
